{"js": "// The title paragraph (\"Quarto Basics\", style \"Title\") is rewritten to\n// \"Thesis on Open Sciene and data - Title WIP\", reproducing the same\n// one-run-per-word/space run layout seen in the target OOXML diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the Title-styled paragraph (the document's first paragraph).\nlet titlePara = null;\nfor (const p of paragraphs.items) {\n  p.load(\"style\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.style === \"Title\") {\n    titlePara = p;\n    break;\n  }\n}\nif (!titlePara) {\n  titlePara = paragraphs.items[0];\n}\n\n// Tokens of the new title, split exactly like the diff's run list\n// (one run per word and one run per separating space).\nconst tokens = [\n  \"Thesis\", \" \", \"on\", \" \", \"Open\", \" \", \"Sciene\", \" \", \"and\", \" \",\n  \"data\", \" \", \"-\", \" \", \"Title\", \" \", \"WIP\"\n];\nconst runsXml = tokens\n  .map((t) => `<w:r><w:t xml:space=\"preserve\">${t}</w:t></w:r>`)\n  .join(\"\");\n\n// Build a minimal OOXML package fragment that replaces the paragraph's\n// range (pPr + runs) while keeping its original \"Title\" style.\nconst ooxml =\n  `<?xml version=\"1.0\" standalone=\"yes\"?>` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">` +\n  `<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">` +\n  `<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>` +\n  `</Relationships></pkg:xmlData></pkg:part>` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body><w:p><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>${runsXml}</w:p></w:body>` +\n  `</w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nconst range = titlePara.getRange();\nrange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The title paragraph (\"Quarto Basics\", style \"Title\") is rewritten to\n# \"Thesis on Open Sciene and data - Title WIP\", reproducing the same\n# one-run-per-word/space run layout seen in the target OOXML diff.\n$d = $word.ActiveDocument\n\n# Locate the Title-styled paragraph (the document's first paragraph).\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n        break\n    }\n}\nif ($titlePara -eq $null) {\n    $titlePara = $d.Paragraphs.Item(1)\n}\n\n# Tokens of the new title, split exactly like the diff's run list\n# (one run per word and one run per separating space).\n$tokens = @(\"Thesis\",\" \",\"on\",\" \",\"Open\",\" \",\"Sciene\",\" \",\"and\",\" \",\"data\",\" \",\"-\",\" \",\"Title\",\" \",\"WIP\")\n$runsXml = \"\"\nforeach ($t in $tokens) {\n    $runsXml += '<w:r><w:t xml:space=\"preserve\">' + $t + '</w:t></w:r>'\n}\n\n# Replace the paragraph's range (pPr + runs) while keeping its\n# original \"Title\" style, via a raw OOXML fragment.\n$payload = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>' + $runsXml + '</w:p>'\n$titlePara.Range.InsertXML($payload)\n"}
